# Re-bucket the age clusters from 13-26/27-40/41-53/54-64 to
# 13-24/25-38/39-52/53-64 across all four sheets, and refresh the
# Revenue / MAU aggregates that were recomputed against the new buckets.

$wb = $excel.ActiveWorkbook

$newHeaders = @("13-24", "25-38", "39-52", "53-64")

# --- DAU: header text only (underlying daily data is unchanged) ---
$wsDau = $wb.Worksheets.Item("DAU")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $wsDau.Cells.Item(1, 2 + $i).Value = $newHeaders[$i]
}

# --- New Users: header text only (underlying daily data is unchanged) ---
$wsNewUsers = $wb.Worksheets.Item("New Users")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $wsNewUsers.Cells.Item(1, 2 + $i).Value = $newHeaders[$i]
}

# --- Revenue: header text + recomputed daily values for the new buckets ---
$wsRevenue = $wb.Worksheets.Item("Revenue")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $wsRevenue.Cells.Item(1, 2 + $i).Value = $newHeaders[$i]
}

$revenueData = @(
    @(1814,1988,2217,2139),
    @(3858,4634,4614,3353),
    @(5371,6149,6108,5223),
    @(8003,8159,8068,8056),
    @(9475,11190,10526,9136),
    @(11365,14003,11535,11755),
    @(11545,13527,14331,13041),
    @(13874,17162,17332,15437),
    @(16245,19649,19033,14716),
    @(19136,20676,21334,17253),
    @(21625,23426,22869,21653),
    @(22608,26218,25077,21352),
    @(22981,25777,27485,23728),
    @(24850,28469,29964,26478),
    @(27991,30257,31238,27022),
    @(28260,32795,34141,29348),
    @(33646,35554,37675,29722),
    @(34595,40595,37125,32944),
    @(34141,40661,42326,36148),
    @(36586,44057,42470,37292),
    @(40775,45303,44598,38076),
    @(40885,44304,49117,40823),
    @(41177,50020,46553,42728),
    @(45135,51127,53057,44227),
    @(45115,53853,52998,46535),
    @(47990,55380,58141,47137),
    @(49393,56084,59512,49549),
    @(51719,62125,61858,52263),
    @(53491,62242,62005,53021),
    @(56501,65772,65187,59667)
)

for ($i = 0; $i -lt $revenueData.Length; $i++) {
    $row = 2 + $i
    $vals = $revenueData[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $wsRevenue.Cells.Item($row, 2 + $j).Value = $vals[$j]
    }
}

# --- MAU: header text + recomputed monthly totals for the new buckets ---
$wsMau = $wb.Worksheets.Item("MAU")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $wsMau.Cells.Item(1, 1 + $i).Value = $newHeaders[$i]
}

$mauData = @(818063, 953926, 951615, 816256)
for ($i = 0; $i -lt $mauData.Length; $i++) {
    $wsMau.Cells.Item(2, 1 + $i).Value = $mauData[$i]
}
